$wb = $excel.ActiveWorkbook

# --- bettaMax sheet ("Store" data): add STORE_NAME column ---
$wsMax = $wb.Worksheets.Item("bettaMax")
$wsMax.Range("D1").Value = "STORE_NAME"
$wsMax.Range("D2").Value = "Hien store"

# Column D width (~17.3 chars); engine snaps width to whole-pixel grid so
# this is the closest achievable value to the source 17.296875.
$wsMax.Columns.Item(4).ColumnWidth = 16.5

# --- Make "bettaMax" the active sheet / selection, as in the source edit ---
$wsMax.Activate()
$wsMax.Range("J25").Select()
